$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.285523653030396
$ws.Range("B1").Value = 2.928736686706543
$ws.Range("C1").Value = 1.878665089607239
$ws.Range("D1").Value = 1.608091711997986
$ws.Range("E1").Value = 1.551949143409729
